# Pavel - new user for linking test
# Adds a new "Linking_AutoUser" row to the Users sheet (row 52), matching
# the format of the existing data rows, then leaves the new row selected
# (mirrors what happens when a user types a new row at the bottom of the
# table and Excel saves with that row/range selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Duplicate the formatting of an existing, fully-populated data row (row 38
# has text in every column, including G, with no hyperlink) onto the new
# row 52, then overwrite the values for the new user.
$ws.Range("A38:G38").Copy()
$ws.Range("A52:G52").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A52").Value = "Linking_AutoUser"
$ws.Range("B52").Value = "Password1"
$ws.Range("C52").ClearContents()
$ws.Range("D52").ClearContents()
$ws.Range("E52").Value = "Default user for Linking tests"
$ws.Range("F52").Value = "N"
$ws.Range("G52").Value = "linking.autouser@mailinator.com"

# Scroll/select like Excel would after entering data on row 52.
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("A52:G52").Select()
